$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add row 8: index 6, period 25/07, rendimento 1.077.343,87
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A8").Value = 6

$ws.Range("B8").Value = "25/07"
$ws.Range("C8").Value = "1.077.343,87"

# Add row 9: index 7, period 25/08, rendimento 956.265,43
$ws.Range("A7").Copy()
$ws.Range("A9").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A9").Value = 7

$ws.Range("B9").Value = "25/08"
$ws.Range("C9").Value = "956.265,43"

$excel.CutCopyMode = $false
